$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 45.89896166666667
$ws.Range("H2").Value = 137.696885
$ws.Range("I2").Value = 0.5698328702801992
$ws.Range("J2").Value = 0.5698328702801992
$ws.Range("M2").Value = 2.759544333333333
$ws.Range("N2").Value = 8.278632999999999
$ws.Range("O2").Value = 0.2574067337278401
$ws.Range("P2").Value = 0.2574067337278401
$ws.Range("Q2").Value = 126.6602195731339
$ws.Range("R2").Value = 1139.941976158205
$ws.Range("S2").Value = 0.1466788179095861
$ws.Range("T2").Value = 0.1466788179095861
$ws.Range("G3").Value = 45.89896166666667
$ws.Range("H3").Value = 137.696885
$ws.Range("I3").Value = 0.5698328702801992
$ws.Range("J3").Value = 0.5698328702801992
$ws.Range("O3").Value = 0.6758254232987829
$ws.Range("P3").Value = 0.6758254232987829
$ws.Range("Q3").Value = 332.5483963392984
$ws.Range("R3").Value = 2992.935567053685
$ws.Range("S3").Value = 0.3851075407666761
$ws.Range("T3").Value = 0.3851075407666761
$ws.Range("G4").Value = 45.89896166666667
$ws.Range("H4").Value = 137.696885
$ws.Range("I4").Value = 0.5698328702801992
$ws.Range("J4").Value = 0.5698328702801992
$ws.Range("M4").Value = 0.5200313333333334
$ws.Range("N4").Value = 1.560094
$ws.Range("O4").Value = 0.0485078515798926
$ws.Range("P4").Value = 0.0485078515798926
$ws.Range("Q4").Value = 23.86889823413223
$ws.Range("R4").Value = 214.82008410719
$ws.Range("S4").Value = 0.0276413682968961
$ws.Range("T4").Value = 0.0276413682968961
$ws.Range("G5").Value = 45.89896166666667
$ws.Range("H5").Value = 137.696885
$ws.Range("I5").Value = 0.5698328702801992
$ws.Range("J5").Value = 0.5698328702801992
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1957573333333333
$ws.Range("N5").Value = 0.587272
$ws.Range("O5").Value = 0.01825999139348442
$ws.Range("P5").Value = 0.01825999139348442
$ws.Range("Q5").Value = 8.985058338635557
$ws.Range("R5").Value = 80.86552504772001
$ws.Range("S5").Value = 0.01040514330704096
$ws.Range("T5").Value = 0.01040514330704096
$ws.Range("I6").Value = 0.04736372570041834
$ws.Range("J6").Value = 0.04736372570041834
$ws.Range("M6").Value = 2.759544333333333
$ws.Range("N6").Value = 8.278632999999999
$ws.Range("O6").Value = 0.2574067337278401
$ws.Range("P6").Value = 0.2574067337278401
$ws.Range("Q6").Value = 10.52782352493422
$ws.Range("R6").Value = 94.750411724408
$ws.Range("S6").Value = 0.01219174192972604
$ws.Range("T6").Value = 0.01219174192972604
$ws.Range("I7").Value = 0.04736372570041834
$ws.Range("J7").Value = 0.04736372570041834
$ws.Range("O7").Value = 0.6758254232987829
$ws.Range("P7").Value = 0.6758254232987829
$ws.Range("S7").Value = 0.03200960997049267
$ws.Range("T7").Value = 0.03200960997049267
$ws.Range("I8").Value = 0.04736372570041834
$ws.Range("J8").Value = 0.04736372570041834
$ws.Range("M8").Value = 0.5200313333333334
$ws.Range("N8").Value = 1.560094
$ws.Range("O8").Value = 0.0485078515798926
$ws.Range("P8").Value = 0.0485078515798926
$ws.Range("Q8").Value = 1.983950045171556
$ws.Range("R8").Value = 17.855550406544
$ws.Range("S8").Value = 0.002297512576546638
$ws.Range("T8").Value = 0.002297512576546638
$ws.Range("I9").Value = 0.04736372570041834
$ws.Range("J9").Value = 0.04736372570041834
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.1957573333333333
$ws.Range("N9").Value = 0.587272
$ws.Range("O9").Value = 0.01825999139348442
$ws.Range("P9").Value = 0.01825999139348442
$ws.Range("Q9").Value = 0.7468257110968889
$ws.Range("R9").Value = 6.721431399872
$ws.Range("S9").Value = 0.0008648612236529958
$ws.Range("T9").Value = 0.0008648612236529958
$ws.Range("G10").Value = 4.651706333333334
$ws.Range("H10").Value = 13.955119
$ws.Range("I10").Value = 0.05775065655894644
$ws.Range("J10").Value = 0.05775065655894644
$ws.Range("M10").Value = 2.759544333333333
$ws.Range("N10").Value = 8.278632999999999
$ws.Range("O10").Value = 0.2574067337278401
$ws.Range("P10").Value = 0.2574067337278401
$ws.Range("Q10").Value = 12.83658985248078
$ws.Range("R10").Value = 115.529308672327
$ws.Range("S10").Value = 0.01486540787547667
$ws.Range("T10").Value = 0.01486540787547667
$ws.Range("G11").Value = 4.651706333333334
$ws.Range("H11").Value = 13.955119
$ws.Range("I11").Value = 0.05775065655894644
$ws.Range("J11").Value = 0.05775065655894644
$ws.Range("O11").Value = 0.6758254232987829
$ws.Range("P11").Value = 0.6758254232987829
$ws.Range("Q11").Value = 33.70266832233767
$ws.Range("R11").Value = 303.324014901039
$ws.Range("S11").Value = 0.03902936191473261
$ws.Range("T11").Value = 0.03902936191473261
$ws.Range("G12").Value = 4.651706333333334
$ws.Range("H12").Value = 13.955119
$ws.Range("I12").Value = 0.05775065655894644
$ws.Range("J12").Value = 0.05775065655894644
$ws.Range("M12").Value = 0.5200313333333334
$ws.Range("N12").Value = 1.560094
$ws.Range("O12").Value = 0.0485078515798926
$ws.Range("P12").Value = 0.0485078515798926
$ws.Range("Q12").Value = 2.419033046798445
$ws.Range("R12").Value = 21.771297421186
$ws.Range("S12").Value = 0.002801360277002725
$ws.Range("T12").Value = 0.002801360277002725
$ws.Range("G13").Value = 4.651706333333334
$ws.Range("H13").Value = 13.955119
$ws.Range("I13").Value = 0.05775065655894644
$ws.Range("J13").Value = 0.05775065655894644
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.1957573333333333
$ws.Range("N13").Value = 0.587272
$ws.Range("O13").Value = 0.01825999139348442
$ws.Range("P13").Value = 0.01825999139348442
$ws.Range("Q13").Value = 0.9106056272631113
$ws.Range("R13").Value = 8.195450645368002
$ws.Range("S13").Value = 0.001054526491734437
$ws.Range("T13").Value = 0.001054526491734437
$ws.Range("G14").Value = 26.182385
$ws.Range("H14").Value = 78.547155
$ws.Range("I14").Value = 0.325052747460436
$ws.Range("J14").Value = 0.325052747460436
$ws.Range("M14").Value = 2.759544333333333
$ws.Range("N14").Value = 8.278632999999999
$ws.Range("O14").Value = 0.2574067337278401
$ws.Range("P14").Value = 0.2574067337278401
$ws.Range("Q14").Value = 72.25145215990166
$ws.Range("R14").Value = 650.2630694391149
$ws.Range("S14").Value = 0.08367076601305129
$ws.Range("T14").Value = 0.08367076601305129
$ws.Range("G15").Value = 26.182385
$ws.Range("H15").Value = 78.547155
$ws.Range("I15").Value = 0.325052747460436
$ws.Range("J15").Value = 0.325052747460436
$ws.Range("O15").Value = 0.6758254232987829
$ws.Range("P15").Value = 0.6758254232987829
$ws.Range("Q15").Value = 189.697322726395
$ws.Range("R15").Value = 1707.275904537555
$ws.Range("S15").Value = 0.2196789106468816
$ws.Range("T15").Value = 0.2196789106468816
$ws.Range("G16").Value = 26.182385
$ws.Range("H16").Value = 78.547155
$ws.Range("I16").Value = 0.325052747460436
$ws.Range("J16").Value = 0.325052747460436
$ws.Range("M16").Value = 0.5200313333333334
$ws.Range("N16").Value = 1.560094
$ws.Range("O16").Value = 0.0485078515798926
$ws.Range("P16").Value = 0.0485078515798926
$ws.Range("Q16").Value = 13.61566058139667
$ws.Range("R16").Value = 122.54094523257
$ws.Range("S16").Value = 0.01576761042944714
$ws.Range("T16").Value = 0.01576761042944714
$ws.Range("G17").Value = 26.182385
$ws.Range("H17").Value = 78.547155
$ws.Range("I17").Value = 0.325052747460436
$ws.Range("J17").Value = 0.325052747460436
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.1957573333333333
$ws.Range("N17").Value = 0.587272
$ws.Range("O17").Value = 0.01825999139348442
$ws.Range("P17").Value = 0.01825999139348442
$ws.Range("Q17").Value = 5.125393867906666
$ws.Range("R17").Value = 46.12854481116
$ws.Range("S17").Value = 0.005935460371056026
$ws.Range("T17").Value = 0.005935460371056026
